$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-12-24 06:30:45"
$newValue = "2025-12-24 06:39:24"

for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
